$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells that change to stay as text,
# since several look like plain numbers (e.g. "22.60", "0.07690")
# and would otherwise be auto-converted/truncated by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.161.47'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").Value = '1.800.26'
$ws.Range("E3").Value = '  +2.43%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '337.86'
$ws.Range("E5").Value = '  +1.01%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '0.4749'
$ws.Range("E7").Value = '  +25.68%  '
$ws.Range("D8").Value = '0.3743'
$ws.Range("E8").Value = '  +11.44%  '
$ws.Range("D9").Value = '45.56'
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").Value = '0.07690'
$ws.Range("E10").Value = '  +6.72%  '
$ws.Range("D11").Value = '1.149'
$ws.Range("E11").Value = '  +1.84%  '
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '22.60'
$ws.Range("E12").Value = '  -0.27%  '
$ws.Range("B13").Value = 'BinanceUSD'
$ws.Range("C13").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D13").Value = '1.002'
$ws.Range("E13").Value = '  +0.01%  '
$ws.Range("D14").Value = '6.393'
$ws.Range("E14").Value = '  +3.62%  '
$ws.Range("D15").Value = '7.450'
$ws.Range("E15").Value = '  +3.52%  '
$ws.Range("D16").Value = '1.797.10'
$ws.Range("E16").Value = '  +2.18%  '
$ws.Range("D17").Value = '0.00001096'
$ws.Range("E17").Value = '  +3.98%  '
$ws.Range("D18").Value = '0.06737'
$ws.Range("E18").Value = '  +2.39%  '
$ws.Range("D19").Value = '82.73'
$ws.Range("E19").Value = '  +2.33%  '
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("D21").Value = '17.55'
$ws.Range("E21").Value = '  +3.19%  '
$ws.Range("D22").Value = '6.448'
$ws.Range("E22").Value = '  +2.84%  '
$ws.Range("D23").Value = '28.161.46'
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("E24").Value = '  +2.78%  '
$ws.Range("D25").Value = '2.408'
$ws.Range("E25").Value = '  +0.51%  '
$ws.Range("D26").Value = '20.89'
$ws.Range("E26").Value = '  +4.94%  '
$ws.Range("D27").Value = '2.411'
$ws.Range("E27").Value = '  +3.35%  '
$ws.Range("D28").Value = '151.68'
$ws.Range("E28").Value = '  -1.11%  '
$ws.Range("D29").Value = '2.002.99'
$ws.Range("E29").Value = '  +2.21%  '
$ws.Range("D30").Value = '134.28'
$ws.Range("E30").Value = '  +1.67%  '
$ws.Range("D31").Value = '1.269'
$ws.Range("E31").Value = '  +0.61%  '
$ws.Range("D32").Value = '4.052'
$ws.Range("E32").Value = '  +0.69%  '
$ws.Range("D33").Value = '0.09676'
$ws.Range("E33").Value = '  +10.55%  '
$ws.Range("D34").Value = '5.966'
$ws.Range("E34").Value = '  +2.61%  '
$ws.Range("E35").Value = '  +2.47%  '
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").Value = '12.26'
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").Value = '0.2229'
$ws.Range("E37").Value = '  +5.41%  '
$ws.Range("D38").Value = '0.06394'
$ws.Range("E38").Value = '  +2.95%  '
$ws.Range("D39").Value = '0.6727'
$ws.Range("E39").Value = '  +0.91%  '
$ws.Range("D40").Value = '5.281'
$ws.Range("E40").Value = '  +2.22%  '
$ws.Range("D41").Value = '1.236'
$ws.Range("E41").Value = '  +1.53%  '
$ws.Range("D42").Value = '1.485'
$ws.Range("E42").Value = '  +2.53%  '
$ws.Range("D43").Value = '8.099'
$ws.Range("E43").Value = '  +1.00%  '
$ws.Range("D44").Value = '14.26'
$ws.Range("E44").Value = '  +3.57%  '
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("D46").Value = '0.6180'
$ws.Range("E46").Value = '  +2.09%  '
$ws.Range("D47").Value = '3.866'
$ws.Range("E47").Value = '  +1.43%  '
$ws.Range("D48").Value = '130.23'
$ws.Range("E48").Value = '  +0.58%  '
$ws.Range("D49").Value = '2.074'
$ws.Range("E49").Value = '  +2.93%  '
$ws.Range("D50").Value = '1.185'
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("D51").Value = '0.07120'
$ws.Range("E51").Value = '  -1.34%  '
